$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 59 on ALC
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

# Row 101 on ALC
$ws.Range("H101").Value = 2000408.2
$ws.Range("I101").Value = 2500464.2
$ws.Range("J101").Value = 184
$ws.Range("K101").Value = 7501392.600000001
$ws.Range("L101").Value = 552
$ws.Range("M101").Value = -7499770.600000001
$ws.Range("N101").Value = -3796

# Row 121 on ALC
$ws.Range("H121").Value = 67139.47
$ws.Range("J121").Value = 67139.47
$ws.Range("L121").Value = 201418.41
$ws.Range("N121").Value = -204912.41

# Row 132 on ALC
$ws.Range("H132").Value = 1623.1034
$ws.Range("I132").Value = 1502.5
$ws.Range("K132").Value = 4507.5
$ws.Range("M132").Value = -1977.5

# Row 138 on ALC
$ws.Range("H138").Value = 2118.394
$ws.Range("I138").Value = 1358.9048
$ws.Range("K138").Value = 4076.7144
$ws.Range("M138").Value = 1063.2856

$ws = $wb.Worksheets.Item("ARM")
# Row 32 on ARM
$ws.Range("H32").Value = 3645117.2
$ws.Range("I32").Value = 719073.8
$ws.Range("K32").Value = 719073.8
$ws.Range("M32").Value = -718786.8

# Row 57 on ARM
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4516

# Row 74 on ARM
$ws.Range("H74").Value = 2269.111
$ws.Range("I74").Value = 1353.4762
$ws.Range("J74").Value = 5473.8335
$ws.Range("K74").Value = 1353.4762
$ws.Range("L74").Value = 5473.8335
$ws.Range("M74").Value = -479.4762000000001
$ws.Range("N74").Value = -7221.8335

# Row 77 on ARM
$ws.Range("H77").Value = 2269.111
$ws.Range("I77").Value = 1353.4762
$ws.Range("J77").Value = 5473.8335
$ws.Range("K77").Value = 6767.381
$ws.Range("L77").Value = 27369.1675
$ws.Range("M77").Value = -2399.381
$ws.Range("N77").Value = -36105.1675

# Row 110 on ARM
$ws.Range("H110").Value = 3421.75
$ws.Range("I110").Value = 2562.3333
$ws.Range("K110").Value = 2562.3333
$ws.Range("M110").Value = -517.3332999999998

# Row 122 on ARM
$ws.Range("H122").Value = 2924.24
$ws.Range("I122").Value = 2005.0555
$ws.Range("K122").Value = 6015.166499999999
$ws.Range("M122").Value = -3565.166499999999

# Row 132 on ARM
$ws.Range("H132").Value = 3707.2964
$ws.Range("I132").Value = 2920.8333
$ws.Range("K132").Value = 8762.499899999999
$ws.Range("M132").Value = -6232.499899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 9 on BSM
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

# Row 107 on BSM
$ws.Range("H107").Value = 2538.282
$ws.Range("I107").Value = 2471.8928
$ws.Range("K107").Value = 2471.8928
$ws.Range("M107").Value = -551.8928000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31 on CRP
$ws.Range("H31").Value = 4372.3335
$ws.Range("I31").Value = 1372.4
$ws.Range("J31").Value = 8122.25
$ws.Range("K31").Value = 1372.4
$ws.Range("L31").Value = 8122.25
$ws.Range("M31").Value = -1077.4
$ws.Range("N31").Value = -8712.25

# Row 34 on CRP
$ws.Range("H34").Value = 4372.3335
$ws.Range("I34").Value = 1372.4
$ws.Range("J34").Value = 8122.25
$ws.Range("K34").Value = 1372.4
$ws.Range("L34").Value = 8122.25
$ws.Range("M34").Value = -1170.4
$ws.Range("N34").Value = -8526.25

# Row 58 on CRP
$ws.Range("H58").Value = 1021.6
$ws.Range("I58").Value = 1103.6666
$ws.Range("J58").Value = 898.5
$ws.Range("K58").Value = 1103.6666
$ws.Range("L58").Value = 898.5
$ws.Range("M58").Value = -900.6666
$ws.Range("N58").Value = -1304.5

# Row 107 on CRP
$ws.Range("H107").Value = 984.9
$ws.Range("I107").Value = 1141.7333
$ws.Range("J107").Value = 514.4
$ws.Range("K107").Value = 1141.7333
$ws.Range("L107").Value = 514.4
$ws.Range("M107").Value = 778.2666999999999
$ws.Range("N107").Value = -4354.4

# Row 132 on CRP
$ws.Range("H132").Value = 2146.0488
$ws.Range("I132").Value = 2092.9412
$ws.Range("J132").Value = 2404
$ws.Range("K132").Value = 6278.823600000001
$ws.Range("L132").Value = 7212
$ws.Range("M132").Value = -3748.823600000001
$ws.Range("N132").Value = -12272

# Row 134 on CRP
$ws.Range("H134").Value = 1277.7693
$ws.Range("I134").Value = 1198.7142
$ws.Range("J134").Value = 1370
$ws.Range("K134").Value = 3596.1426
$ws.Range("L134").Value = 4110
$ws.Range("M134").Value = -1061.1426
$ws.Range("N134").Value = -9180

# Row 136 on CRP
$ws.Range("H136").Value = 1021.6
$ws.Range("I136").Value = 1103.6666
$ws.Range("J136").Value = 898.5
$ws.Range("K136").Value = 3310.9998
$ws.Range("L136").Value = 2695.5
$ws.Range("M136").Value = -760.9998000000001
$ws.Range("N136").Value = -7795.5

$ws = $wb.Worksheets.Item("CUL")
# Row 107 on CUL
$ws.Range("H107").Value = 703.1539
$ws.Range("J107").Value = 706.7143
$ws.Range("L107").Value = 2120.1429
$ws.Range("N107").Value = -5960.1429

# Row 132 on CUL
$ws.Range("H132").Value = 1391.2858
$ws.Range("I132").Value = 1110.5714
$ws.Range("J132").Value = 1484.8572
$ws.Range("K132").Value = 9995.142600000001
$ws.Range("L132").Value = 13363.7148
$ws.Range("M132").Value = -7465.142600000001
$ws.Range("N132").Value = -18423.7148

$ws = $wb.Worksheets.Item("GSM")
# Row 96 on GSM
$ws.Range("H96").Value = 30261
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 30261
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 30261
$ws.Range("N96").Value = -35753
$ws.Range("M96").ClearContents()

# Row 99 on GSM
$ws.Range("H99").Value = 7442.091
$ws.Range("I99").Value = 2371.6667
$ws.Range("J99").Value = 30259
$ws.Range("K99").Value = 2371.6667
$ws.Range("L99").Value = 30259
$ws.Range("M99").Value = -125.6667000000002
$ws.Range("N99").Value = -34751

# Row 122 on GSM
$ws.Range("H122").Value = 3329.9092
$ws.Range("I122").Value = 1322.2
$ws.Range("J122").Value = 5003
$ws.Range("K122").Value = 3966.6
$ws.Range("L122").Value = 15009
$ws.Range("M122").Value = -1516.6
$ws.Range("N122").Value = -19909

# Row 126 on GSM
$ws.Range("H126").Value = 51388.77
$ws.Range("I126").Value = 2818.8333
$ws.Range("J126").Value = 93020.14
$ws.Range("K126").Value = 8456.499899999999
$ws.Range("L126").Value = 279060.42
$ws.Range("M126").Value = -5986.499899999999
$ws.Range("N126").Value = -284000.42

$ws = $wb.Worksheets.Item("LTW")
# Row 16 on LTW
$ws.Range("H16").Value = 1051.5883
$ws.Range("I16").Value = 978.4666999999999
$ws.Range("K16").Value = 978.4666999999999
$ws.Range("M16").Value = -808.4666999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107 on WVR
$ws.Range("H107").Value = 45508120
$ws.Range("I107").Value = 2405.889
$ws.Range("K107").Value = 7217.667
$ws.Range("M107").Value = -5297.667

# Row 132 on WVR
$ws.Range("H132").Value = 1102036.4
$ws.Range("I132").Value = 1367242.8
$ws.Range("K132").Value = 4101728.4
$ws.Range("M132").Value = -4099198.4

# Row 136 on WVR
$ws.Range("H136").Value = 6337.4287
$ws.Range("I136").Value = 4675.607
$ws.Range("K136").Value = 14026.821
$ws.Range("M136").Value = -11476.821
